$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: column B now holds the "EI(g_wh)" label; the old "old" and
# --- "step_change" columns (C, D) lose their header text entirely.
$ws.Range("B1").Value = "EI(g_wh)"
$ws.Range("C1:D1").ClearContents()

# --- Column B (the surviving EI series) now displays three decimal places.
$ws.Range("B2:B31").NumberFormat = "0.000"

# --- Columns C ("old") and D ("step_change") are no longer populated.
$ws.Range("C2:D31").ClearContents()

# --- The old "interpolated" highlight in column C (and the trailing blank
# --- cell's highlight) is removed now that the column is unused, and a
# --- handful of further blank (but touched/formatted) rows were appended
# --- below the table.
$ws.Range("C20:C36").Interior.ColorIndex = -4142

# --- Restore the selection left behind by the editor.
$ws.Range("K14:K15").Select()
